# Re-apply the table style that PowerPoint assigns to the three balance-sheet
# tables (slides 14, 15, 16) after the "Integral" theme's default table style
# GUID was swapped in for the original one.
#
# Original tableStyleId : {1C41FA46-E823-4A44-9516-01A8700AFEE3}
# New      tableStyleId : {D23BB660-F28B-47AA-AC24-F503BF0C51E2}

$p = $ppt.ActivePresentation

$newStyleId = "{D23BB660-F28B-47AA-AC24-F503BF0C51E2}"
$targetSlides = @(14, 15, 16)

foreach ($slideIdx in $targetSlides) {
    $slide = $p.Slides.Item($slideIdx)

    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}
